$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 9032.571
$ws.Range("I2").Value2 = 1558
$ws.Range("J2").Value2 = 18998.666
$ws.Range("K2").Value2 = 1558
$ws.Range("L2").Value2 = 18998.666
$ws.Range("M2").Value2 = -1445
$ws.Range("N2").Value2 = -19224.666
$ws.Range("H11").Value2 = 4336.2573
$ws.Range("I11").Value2 = 4336.2573
$ws.Range("K11").Value2 = 4336.2573
$ws.Range("M11").Value2 = -4196.2573
$ws.Range("H33").Value2 = 5263650
$ws.Range("I33").Value2 = 8333935
$ws.Range("K33").Value2 = 8333935
$ws.Range("M33").Value2 = -8333706
$ws.Range("H38").Value2 = 444.73334
$ws.Range("I38").Value2 = 155.07143
$ws.Range("K38").Value2 = 465.21429
$ws.Range("M38").Value2 = -93.21429
$ws.Range("H39").Value2 = 455.8125
$ws.Range("I39").Value2 = 235.72728
$ws.Range("J39").Value2 = 940
$ws.Range("K39").Value2 = 707.18184
$ws.Range("L39").Value2 = 2820
$ws.Range("M39").Value2 = -411.18184
$ws.Range("N39").Value2 = -3412
$ws.Range("H41").Value2 = 721.1429
$ws.Range("I41").Value2 = 79.8
$ws.Range("J41").Value2 = 2324.5
$ws.Range("K41").Value2 = 79.8
$ws.Range("L41").Value2 = 2324.5
$ws.Range("M41").Value2 = 360.2
$ws.Range("N41").Value2 = -3204.5
$ws.Range("H43").Value2 = 3121.75
$ws.Range("I43").Value2 = 1927.3334
$ws.Range("J43").Value2 = 3838.4
$ws.Range("K43").Value2 = 1927.3334
$ws.Range("L43").Value2 = 3838.4
$ws.Range("M43").Value2 = -1858.3334
$ws.Range("N43").Value2 = -3976.4
$ws.Range("H62").Value2 = 44446976
$ws.Range("I62").Value2 = 44446976
$ws.Range("K62").Value2 = 44446976
$ws.Range("M62").Value2 = -44446352
$ws.Range("H65").Value2 = 44446976
$ws.Range("I65").Value2 = 44446976
$ws.Range("K65").Value2 = 222234880
$ws.Range("M65").Value2 = -222231760
$ws.Range("H70").Value2 = 5530.263
$ws.Range("I70").Value2 = 4814.4287
$ws.Range("J70").Value2 = 5947.8335
$ws.Range("K70").Value2 = 14443.2861
$ws.Range("L70").Value2 = 17843.5005
$ws.Range("M70").Value2 = -14173.2861
$ws.Range("N70").Value2 = -18383.5005
$ws.Range("H73").Value2 = 5530.263
$ws.Range("I73").Value2 = 4814.4287
$ws.Range("J73").Value2 = 5947.8335
$ws.Range("K73").Value2 = 14443.2861
$ws.Range("L73").Value2 = 17843.5005
$ws.Range("M73").Value2 = -13507.2861
$ws.Range("N73").Value2 = -19715.5005
$ws.Range("H74").Value2 = 6666.467
$ws.Range("I74").Value2 = 4999.7
$ws.Range("J74").Value2 = 10000
$ws.Range("K74").Value2 = 4999.7
$ws.Range("L74").Value2 = 10000
$ws.Range("M74").Value2 = -4063.7
$ws.Range("N74").Value2 = -11872
$ws.Range("H77").Value2 = 6666.467
$ws.Range("I77").Value2 = 4999.7
$ws.Range("J77").Value2 = 10000
$ws.Range("K77").Value2 = 24998.5
$ws.Range("L77").Value2 = 50000
$ws.Range("M77").Value2 = -20318.5
$ws.Range("N77").Value2 = -59360
$ws.Range("H88").Value2 = 13140.8
$ws.Range("J88").Value2 = 11801
$ws.Range("L88").Value2 = 11801
$ws.Range("N88").Value2 = -12613
$ws.Range("H91").Value2 = 13140.8
$ws.Range("J91").Value2 = 11801
$ws.Range("L91").Value2 = 11801
$ws.Range("N91").Value2 = -14609
$ws.Range("H92").Value2 = 1955
$ws.Range("I92").Value2 = 1909.6154
$ws.Range("K92").Value2 = 1909.6154
$ws.Range("M92").Value2 = -661.6153999999999
$ws.Range("H94").Value2 = 1999
$ws.Range("I94").Value2 = 1999
$ws.Range("K94").Value2 = 1999
$ws.Range("M94").Value2 = -1548
$ws.Range("H96").Value2 = 9098934
$ws.Range("I96").Value2 = 330.75
$ws.Range("J96").Value2 = 14298136
$ws.Range("K96").Value2 = 992.25
$ws.Range("L96").Value2 = 42894408
$ws.Range("M96").Value2 = 380.75
$ws.Range("N96").Value2 = -42897154
$ws.Range("H98").Value2 = 3451.182
$ws.Range("J98").Value2 = 3064
$ws.Range("L98").Value2 = 3064
$ws.Range("N98").Value2 = -6060
$ws.Range("H101").Value2 = 2469.8572
$ws.Range("I101").Value2 = 1663.8
$ws.Range("K101").Value2 = 4991.4
$ws.Range("M101").Value2 = -3369.4
$ws.Range("H106").Value2 = 1014.9048
$ws.Range("I106").Value2 = 1014.9048
$ws.Range("K106").Value2 = 1014.9048
$ws.Range("M106").Value2 = -383.9048
$ws.Range("H107").Value2 = 33335744
$ws.Range("I107").Value2 = 38463820
$ws.Range("J107").Value2 = 3250
$ws.Range("K107").Value2 = 38463820
$ws.Range("L107").Value2 = 3250
$ws.Range("M107").Value2 = -38461900
$ws.Range("N107").Value2 = -7090
$ws.Range("H122").Value2 = 3451.182
$ws.Range("J122").Value2 = 3064
$ws.Range("L122").Value2 = 9192
$ws.Range("N122").Value2 = -14092
$ws.Range("H125").Value2 = 2496.6667
$ws.Range("I125").Value2 = 2495
$ws.Range("J125").Value2 = 2497.5
$ws.Range("K125").Value2 = 22455
$ws.Range("L125").Value2 = 22477.5
$ws.Range("M125").Value2 = -19995
$ws.Range("N125").Value2 = -27397.5
$ws.Range("H132").Value2 = 1241.4348
$ws.Range("I132").Value2 = 1073.9524
$ws.Range("K132").Value2 = 3221.857199999999
$ws.Range("M132").Value2 = -691.8571999999995
$ws.Range("H137").Value2 = 14713581
$ws.Range("I137").Value2 = 27778960
$ws.Range("J137").Value2 = 15028.125
$ws.Range("K137").Value2 = 83336880
$ws.Range("L137").Value2 = 45084.375
$ws.Range("M137").Value2 = -83334330
$ws.Range("N137").Value2 = -50184.375
$ws.Range("H138").Value2 = 3787.9688
$ws.Range("J138").Value2 = 3785
$ws.Range("L138").Value2 = 11355
$ws.Range("N138").Value2 = -21635
$ws.Range("H141").Value2 = 3284.5
$ws.Range("I141").Value2 = 1569
$ws.Range("K141").Value2 = 4707
$ws.Range("M141").Value2 = 473

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 18184522
$ws.Range("J2").Value2 = 1010
$ws.Range("L2").Value2 = 1010
$ws.Range("N2").Value2 = -1236
$ws.Range("H28").Value2 = 20892
$ws.Range("I28").Value2 = 20892
$ws.Range("K28").Value2 = 20892
$ws.Range("M28").Value2 = -20700
$ws.Range("H32").Value2 = 62904.277
$ws.Range("I32").Value2 = 62904.277
$ws.Range("K32").Value2 = 62904.277
$ws.Range("M32").Value2 = -62617.277
$ws.Range("H43").Value2 = 35498.6
$ws.Range("J43").Value2 = 34373.75
$ws.Range("L43").Value2 = 34373.75
$ws.Range("N43").Value2 = -34999.75
$ws.Range("H45").Value2 = 2409.875
$ws.Range("I45").Value2 = 1611.2858
$ws.Range("J45").Value2 = 8000
$ws.Range("K45").Value2 = 1611.2858
$ws.Range("L45").Value2 = 8000
$ws.Range("M45").Value2 = -1234.2858
$ws.Range("N45").Value2 = -8754
$ws.Range("H61").Value2 = 4988.724
$ws.Range("I61").Value2 = 3838
$ws.Range("J61").Value2 = 9399.833
$ws.Range("K61").Value2 = 3838
$ws.Range("L61").Value2 = 9399.833
$ws.Range("M61").Value2 = -3626
$ws.Range("N61").Value2 = -9823.833
$ws.Range("H74").Value2 = 3591
$ws.Range("I74").Value2 = 1695.3572
$ws.Range("K74").Value2 = 1695.3572
$ws.Range("M74").Value2 = -821.3571999999999
$ws.Range("H77").Value2 = 3591
$ws.Range("I77").Value2 = 1695.3572
$ws.Range("K77").Value2 = 8476.786
$ws.Range("M77").Value2 = -4108.786
$ws.Range("H80").Value2 = 84999.5
$ws.Range("I80").Value2 = 20000
$ws.Range("J80").Value2 = 149999
$ws.Range("K80").Value2 = 20000
$ws.Range("L80").Value2 = 149999
$ws.Range("M80").Value2 = -19002
$ws.Range("N80").Value2 = -151995
$ws.Range("H83").Value2 = 84999.5
$ws.Range("I83").Value2 = 20000
$ws.Range("J83").Value2 = 149999
$ws.Range("K83").Value2 = 60000
$ws.Range("L83").Value2 = 449997
$ws.Range("M83").Value2 = -55008
$ws.Range("N83").Value2 = -459981
$ws.Range("H99").Value2 = 20892
$ws.Range("I99").Value2 = 20892
$ws.Range("K99").Value2 = 20892
$ws.Range("M99").Value2 = -17897
$ws.Range("H110").Value2 = 154001810
$ws.Range("I110").Value2 = 154001810
$ws.Range("J110").Value2 = 0
$ws.Range("K110").Value2 = 154001810
$ws.Range("L110").Value2 = 0
$ws.Range("M110").Value2 = -153999765
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value2 = 18184522
$ws.Range("J116").Value2 = 1010
$ws.Range("L116").Value2 = 1010
$ws.Range("N116").Value2 = -5598
$ws.Range("H122").Value2 = 2042.2142
$ws.Range("I122").Value2 = 2053.818
$ws.Range("K122").Value2 = 6161.454000000001
$ws.Range("M122").Value2 = -3711.454000000001
$ws.Range("H132").Value2 = 3086.682
$ws.Range("I132").Value2 = 2459.8064
$ws.Range("J132").Value2 = 4581.5386
$ws.Range("K132").Value2 = 7379.4192
$ws.Range("L132").Value2 = 13744.6158
$ws.Range("M132").Value2 = -4849.4192
$ws.Range("N132").Value2 = -18804.6158
$ws.Range("H136").Value2 = 4988.724
$ws.Range("I136").Value2 = 3838
$ws.Range("J136").Value2 = 9399.833
$ws.Range("K136").Value2 = 11514
$ws.Range("L136").Value2 = 28199.499
$ws.Range("M136").Value2 = -8964
$ws.Range("N136").Value2 = -33299.499

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 18184522
$ws.Range("J3").Value2 = 1010
$ws.Range("L3").Value2 = 1010
$ws.Range("N3").Value2 = -1238
$ws.Range("H20").Value2 = 4071.0667
$ws.Range("I20").Value2 = 3962.25
$ws.Range("J20").Value2 = 4195.4287
$ws.Range("K20").Value2 = 3962.25
$ws.Range("L20").Value2 = 4195.4287
$ws.Range("M20").Value2 = -3715.25
$ws.Range("N20").Value2 = -4689.4287
$ws.Range("H86").Value2 = 290253.72
$ws.Range("I86").Value2 = 3825
$ws.Range("J86").Value2 = 576682.44
$ws.Range("K86").Value2 = 3825
$ws.Range("L86").Value2 = 576682.44
$ws.Range("M86").Value2 = -2702
$ws.Range("N86").Value2 = -578928.44
$ws.Range("H89").Value2 = 290253.72
$ws.Range("I89").Value2 = 3825
$ws.Range("J89").Value2 = 576682.44
$ws.Range("K89").Value2 = 19125
$ws.Range("L89").Value2 = 2883412.2
$ws.Range("M89").Value2 = -13509
$ws.Range("N89").Value2 = -2894644.2
$ws.Range("H94").Value2 = 2883.1667
$ws.Range("I94").Value2 = 2574.75
$ws.Range("K94").Value2 = 2574.75
$ws.Range("M94").Value2 = -2123.75
$ws.Range("H99").Value2 = 3231.7856
$ws.Range("I99").Value2 = 2805.8
$ws.Range("K99").Value2 = 2805.8
$ws.Range("M99").Value2 = -1307.8
$ws.Range("H102").Value2 = 18661.4
$ws.Range("I102").Value2 = 8397.429
$ws.Range("J102").Value2 = 42610.668
$ws.Range("K102").Value2 = 8397.429
$ws.Range("L102").Value2 = 42610.668
$ws.Range("M102").Value2 = -5152.429
$ws.Range("N102").Value2 = -49100.668
$ws.Range("H105").Value2 = 111140410
$ws.Range("I105").Value2 = 125032616
$ws.Range("J105").Value2 = 2750
$ws.Range("K105").Value2 = 125032616
$ws.Range("L105").Value2 = 2750
$ws.Range("M105").Value2 = -125030869
$ws.Range("N105").Value2 = -6244
$ws.Range("H107").Value2 = 3990.2307
$ws.Range("I107").Value2 = 3734.0908
$ws.Range("J107").Value2 = 5399
$ws.Range("K107").Value2 = 3734.0908
$ws.Range("L107").Value2 = 5399
$ws.Range("M107").Value2 = -1814.0908
$ws.Range("N107").Value2 = -9239
$ws.Range("H134").Value2 = 2765.7058
$ws.Range("I134").Value2 = 1343.2693
$ws.Range("K134").Value2 = 4029.8079
$ws.Range("M134").Value2 = -1494.8079
$ws.Range("H135").Value2 = 98999
$ws.Range("J135").Value2 = 98999
$ws.Range("L135").Value2 = 98999
$ws.Range("N135").Value2 = -109139

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 215.26666
$ws.Range("I7").Value2 = 84.166664
$ws.Range("K7").Value2 = 84.166664
$ws.Range("M7").Value2 = 28.833336
$ws.Range("H11").Value2 = 2005
$ws.Range("I11").Value2 = 757.5
$ws.Range("J11").Value2 = 4500
$ws.Range("K11").Value2 = 757.5
$ws.Range("L11").Value2 = 4500
$ws.Range("M11").Value2 = -617.5
$ws.Range("N11").Value2 = -4780
$ws.Range("H16").Value2 = 5412.067
$ws.Range("I16").Value2 = 4297.857
$ws.Range("K16").Value2 = 4297.857
$ws.Range("M16").Value2 = -4010.857
$ws.Range("H22").Value2 = 649.7222
$ws.Range("I22").Value2 = 587.3571
$ws.Range("J22").Value2 = 868
$ws.Range("K22").Value2 = 587.3571
$ws.Range("L22").Value2 = 868
$ws.Range("M22").Value2 = -237.3570999999999
$ws.Range("N22").Value2 = -1568
$ws.Range("H41").Value2 = 21348.166
$ws.Range("I41").Value2 = 7331.25
$ws.Range("J41").Value2 = 49382
$ws.Range("K41").Value2 = 7331.25
$ws.Range("L41").Value2 = 49382
$ws.Range("M41").Value2 = -6903.25
$ws.Range("N41").Value2 = -50238
$ws.Range("H69").Value2 = 23515.334
$ws.Range("I69").Value2 = 23515.334
$ws.Range("K69").Value2 = 23515.334
$ws.Range("M69").Value2 = -22766.334
$ws.Range("H70").Value2 = 80090
$ws.Range("J70").Value2 = 80090
$ws.Range("L70").Value2 = 80090
$ws.Range("N70").Value2 = -80720
$ws.Range("H72").Value2 = 23515.334
$ws.Range("I72").Value2 = 23515.334
$ws.Range("K72").Value2 = 70546.002
$ws.Range("M72").Value2 = -66802.002
$ws.Range("H73").Value2 = 80090
$ws.Range("J73").Value2 = 80090
$ws.Range("L73").Value2 = 80090
$ws.Range("N73").Value2 = -82274
$ws.Range("H74").Value2 = 50314
$ws.Range("J74").Value2 = 50314
$ws.Range("L74").Value2 = 50314
$ws.Range("N74").Value2 = -52062
$ws.Range("H75").Value2 = 70000
$ws.Range("J75").Value2 = 70000
$ws.Range("L75").Value2 = 70000
$ws.Range("N75").Value2 = -71996
$ws.Range("H77").Value2 = 50314
$ws.Range("J77").Value2 = 50314
$ws.Range("L77").Value2 = 150942
$ws.Range("N77").Value2 = -159678
$ws.Range("H78").Value2 = 70000
$ws.Range("J78").Value2 = 70000
$ws.Range("L78").Value2 = 210000
$ws.Range("N78").Value2 = -219984
$ws.Range("H82").Value2 = 74950
$ws.Range("J82").Value2 = 74950
$ws.Range("L82").Value2 = 74950
$ws.Range("N82").Value2 = -75672
$ws.Range("H85").Value2 = 74950
$ws.Range("J85").Value2 = 74950
$ws.Range("L85").Value2 = 74950
$ws.Range("N85").Value2 = -77446
$ws.Range("H99").Value2 = 3298.1428
$ws.Range("I99").Value2 = 2515.8333
$ws.Range("J99").Value2 = 7992
$ws.Range("K99").Value2 = 2515.8333
$ws.Range("L99").Value2 = 7992
$ws.Range("M99").Value2 = -1017.8333
$ws.Range("N99").Value2 = -10988
$ws.Range("H113").Value2 = 5412.067
$ws.Range("I113").Value2 = 4297.857
$ws.Range("K113").Value2 = 4297.857
$ws.Range("M113").Value2 = -2127.857
$ws.Range("H120").Value2 = 0
$ws.Range("I120").Value2 = 0
$ws.Range("K120").Value2 = 0
$ws.Range("M120").ClearContents()
$ws.Range("H122").Value2 = 73089.57
$ws.Range("I122").Value2 = 84971.164
$ws.Range("K122").Value2 = 254913.492
$ws.Range("M122").Value2 = -252463.492
$ws.Range("H126").Value2 = 3298.1428
$ws.Range("I126").Value2 = 2515.8333
$ws.Range("J126").Value2 = 7992
$ws.Range("K126").Value2 = 7547.499899999999
$ws.Range("L126").Value2 = 23976
$ws.Range("M126").Value2 = -5077.499899999999
$ws.Range("N126").Value2 = -28916
$ws.Range("H132").Value2 = 8555.091
$ws.Range("J132").Value2 = 8157
$ws.Range("L132").Value2 = 24471
$ws.Range("N132").Value2 = -29531
$ws.Range("H134").Value2 = 9647.866
$ws.Range("I134").Value2 = 9070.182
$ws.Range("K134").Value2 = 27210.546
$ws.Range("M134").Value2 = -24675.546

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 2360.0667
$ws.Range("I2").Value2 = 21.916666
$ws.Range("J2").Value2 = 3918.8333
$ws.Range("K2").Value2 = 131.499996
$ws.Range("L2").Value2 = 23512.9998
$ws.Range("M2").Value2 = -18.49999600000001
$ws.Range("N2").Value2 = -23738.9998
$ws.Range("H7").Value2 = 783.94446
$ws.Range("I7").Value2 = 797.8333
$ws.Range("J7").Value2 = 756.1667
$ws.Range("K7").Value2 = 2393.4999
$ws.Range("L7").Value2 = 2268.5001
$ws.Range("M7").Value2 = -2281.4999
$ws.Range("N7").Value2 = -2492.5001
$ws.Range("H11").Value2 = 7449.875
$ws.Range("I11").Value2 = 9545.583
$ws.Range("K11").Value2 = 28636.749
$ws.Range("M11").Value2 = -28496.749
$ws.Range("H14").Value2 = 803.6
$ws.Range("I14").Value2 = 803.6
$ws.Range("K14").Value2 = 2410.8
$ws.Range("M14").Value2 = -2237.8
$ws.Range("H45").Value2 = 2283
$ws.Range("I45").Value2 = 0
$ws.Range("K45").Value2 = 0
$ws.Range("M45").ClearContents()
$ws.Range("H68").Value2 = 5110.75
$ws.Range("J68").Value2 = 5110.75
$ws.Range("L68").Value2 = 15332.25
$ws.Range("N68").Value2 = -16954.25
$ws.Range("H71").Value2 = 5110.75
$ws.Range("J71").Value2 = 5110.75
$ws.Range("L71").Value2 = 45996.75
$ws.Range("N71").Value2 = -54108.75
$ws.Range("H75").Value2 = 10682.5
$ws.Range("I75").Value2 = 1500
$ws.Range("K75").Value2 = 4500
$ws.Range("M75").Value2 = -3502
$ws.Range("H78").Value2 = 10682.5
$ws.Range("I78").Value2 = 1500
$ws.Range("K78").Value2 = 13500
$ws.Range("M78").Value2 = -8508
$ws.Range("H80").Value2 = 3633.3333
$ws.Range("I80").Value2 = 2800
$ws.Range("J80").Value2 = 4050
$ws.Range("K80").Value2 = 8400
$ws.Range("L80").Value2 = 12150
$ws.Range("M80").Value2 = -7464
$ws.Range("N80").Value2 = -14022
$ws.Range("H81").Value2 = 679.6667
$ws.Range("J81").Value2 = 0
$ws.Range("L81").Value2 = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value2 = 3633.3333
$ws.Range("I83").Value2 = 2800
$ws.Range("J83").Value2 = 4050
$ws.Range("K83").Value2 = 25200
$ws.Range("L83").Value2 = 36450
$ws.Range("M83").Value2 = -20520
$ws.Range("N83").Value2 = -45810
$ws.Range("H84").Value2 = 679.6667
$ws.Range("J84").Value2 = 0
$ws.Range("L84").Value2 = 0
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value2 = 1529.4849
$ws.Range("J107").Value2 = 2150.65
$ws.Range("L107").Value2 = 6451.950000000001
$ws.Range("N107").Value2 = -10291.95
$ws.Range("H110").Value2 = 11010
$ws.Range("I110").Value2 = 6500
$ws.Range("J110").Value2 = 20030
$ws.Range("K110").Value2 = 19500
$ws.Range("L110").Value2 = 60090
$ws.Range("M110").Value2 = -15410
$ws.Range("N110").Value2 = -68270
$ws.Range("H111").Value2 = 17522.5
$ws.Range("I111").Value2 = 10000
$ws.Range("K111").Value2 = 30000
$ws.Range("M111").Value2 = -26933
$ws.Range("H128").Value2 = 443977.4
$ws.Range("I128").Value2 = 443977.4
$ws.Range("K128").Value2 = 1331932.2
$ws.Range("M128").Value2 = -1326952.2
$ws.Range("H129").Value2 = 1357.7858
$ws.Range("I129").Value2 = 693
$ws.Range("J129").Value2 = 10000
$ws.Range("K129").Value2 = 2079
$ws.Range("L129").Value2 = 30000
$ws.Range("M129").Value2 = 2921
$ws.Range("N129").Value2 = -40000
$ws.Range("H132").Value2 = 1438.5555
$ws.Range("I132").Value2 = 1438.5555
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 12946.9995
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -10416.9995
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value2 = 3325.5
$ws.Range("I137").Value2 = 1942.5
$ws.Range("J137").Value2 = 5400
$ws.Range("K137").Value2 = 5827.5
$ws.Range("L137").Value2 = 16200
$ws.Range("M137").Value2 = -727.5
$ws.Range("N137").Value2 = -26400
$ws.Range("H138").Value2 = 5002711
$ws.Range("I138").Value2 = 7144515.5
$ws.Range("J138").Value2 = 5166.3335
$ws.Range("K138").Value2 = 21433546.5
$ws.Range("L138").Value2 = 15499.0005
$ws.Range("M138").Value2 = -21428406.5
$ws.Range("N138").Value2 = -25779.0005
$ws.Range("H139").Value2 = 4697.25
$ws.Range("I139").Value2 = 4013
$ws.Range("K139").Value2 = 12039
$ws.Range("M139").Value2 = -6899
$ws.Range("H140").Value2 = 1315
$ws.Range("I140").Value2 = 1130
$ws.Range("K140").Value2 = 3390
$ws.Range("M140").Value2 = 1790

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value2 = 16526.715
$ws.Range("I43").Value2 = 396
$ws.Range("J43").Value2 = 28624.75
$ws.Range("K43").Value2 = 396
$ws.Range("L43").Value2 = 28624.75
$ws.Range("M43").Value2 = -245
$ws.Range("N43").Value2 = -28926.75
$ws.Range("H53").Value2 = 50605.75
$ws.Range("I53").Value2 = 48342
$ws.Range("J53").Value2 = 51964
$ws.Range("K53").Value2 = 48342
$ws.Range("L53").Value2 = 51964
$ws.Range("M53").Value2 = -47711
$ws.Range("N53").Value2 = -53226
$ws.Range("H63").Value2 = 40001
$ws.Range("J63").Value2 = 40001
$ws.Range("L63").Value2 = 40001
$ws.Range("N63").Value2 = -41373
$ws.Range("H66").Value2 = 40001
$ws.Range("J66").Value2 = 40001
$ws.Range("L66").Value2 = 120003
$ws.Range("N66").Value2 = -126867
$ws.Range("H80").Value2 = 1633.3334
$ws.Range("I80").Value2 = 1450
$ws.Range("J80").Value2 = 2000
$ws.Range("K80").Value2 = 1450
$ws.Range("L80").Value2 = 2000
$ws.Range("M80").Value2 = -452
$ws.Range("N80").Value2 = -3996
$ws.Range("H83").Value2 = 1633.3334
$ws.Range("I83").Value2 = 1450
$ws.Range("J83").Value2 = 2000
$ws.Range("K83").Value2 = 7250
$ws.Range("L83").Value2 = 10000
$ws.Range("M83").Value2 = -2258
$ws.Range("N83").Value2 = -19984
$ws.Range("H93").Value2 = 44999
$ws.Range("J93").Value2 = 44999
$ws.Range("L93").Value2 = 44999
$ws.Range("N93").Value2 = -48743
$ws.Range("H98").Value2 = 50996.668
$ws.Range("J98").Value2 = 50996.668
$ws.Range("L98").Value2 = 50996.668
$ws.Range("N98").Value2 = -56986.668
$ws.Range("H102").Value2 = 1442.2307
$ws.Range("I102").Value2 = 1064.375
$ws.Range("J102").Value2 = 2046.8
$ws.Range("K102").Value2 = 1064.375
$ws.Range("L102").Value2 = 2046.8
$ws.Range("M102").Value2 = 557.625
$ws.Range("N102").Value2 = -5290.8
$ws.Range("H126").Value2 = 3192.9285
$ws.Range("I126").Value2 = 2675.9
$ws.Range("K126").Value2 = 8027.700000000001
$ws.Range("M126").Value2 = -5557.700000000001
$ws.Range("H132").Value2 = 5920.5
$ws.Range("I132").Value2 = 3780
$ws.Range("J132").Value2 = 9488
$ws.Range("K132").Value2 = 11340
$ws.Range("L132").Value2 = 28464
$ws.Range("M132").Value2 = -8810
$ws.Range("N132").Value2 = -33524

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 2548.9
$ws.Range("I16").Value2 = 1855.5714
$ws.Range("K16").Value2 = 1855.5714
$ws.Range("M16").Value2 = -1685.5714
$ws.Range("H22").Value2 = 2291.138
$ws.Range("I22").Value2 = 1382.0286
$ws.Range("K22").Value2 = 1382.0286
$ws.Range("M22").Value2 = -1087.0286
$ws.Range("H27").Value2 = 2291.138
$ws.Range("I27").Value2 = 1382.0286
$ws.Range("K27").Value2 = 1382.0286
$ws.Range("M27").Value2 = -1275.0286
$ws.Range("H46").Value2 = 3332.3416
$ws.Range("J46").Value2 = 4126
$ws.Range("L46").Value2 = 4126
$ws.Range("N46").Value2 = -4502
$ws.Range("H55").Value2 = 258.13635
$ws.Range("I55").Value2 = 230.61111
$ws.Range("J55").Value2 = 382
$ws.Range("K55").Value2 = 230.61111
$ws.Range("L55").Value2 = 382
$ws.Range("M55").Value2 = -57.61111
$ws.Range("N55").Value2 = -728
$ws.Range("H61").Value2 = 3450.889
$ws.Range("I61").Value2 = 3506.8845
$ws.Range("K61").Value2 = 3506.8845
$ws.Range("M61").Value2 = -3304.8845
$ws.Range("H82").Value2 = 2527.5454
$ws.Range("I82").Value2 = 1580.6
$ws.Range("K82").Value2 = 1580.6
$ws.Range("M82").Value2 = -1219.6
$ws.Range("H85").Value2 = 2527.5454
$ws.Range("I85").Value2 = 1580.6
$ws.Range("K85").Value2 = 1580.6
$ws.Range("M85").Value2 = -332.5999999999999
$ws.Range("H93").Value2 = 3099.625
$ws.Range("I93").Value2 = 3099.625
$ws.Range("K93").Value2 = 3099.625
$ws.Range("M93").Value2 = -1851.625
$ws.Range("H97").Value2 = 18465.4
$ws.Range("J97").Value2 = 18465.4
$ws.Range("L97").Value2 = 18465.4
$ws.Range("N97").Value2 = -20447.4
$ws.Range("H113").Value2 = 3450.889
$ws.Range("I113").Value2 = 3506.8845
$ws.Range("K113").Value2 = 3506.8845
$ws.Range("M113").Value2 = -1336.8845
$ws.Range("H122").Value2 = 3777.2222
$ws.Range("I122").Value2 = 3019.8823
$ws.Range("J122").Value2 = 5064.7
$ws.Range("K122").Value2 = 9059.6469
$ws.Range("L122").Value2 = 15194.1
$ws.Range("M122").Value2 = -6609.6469
$ws.Range("N122").Value2 = -20094.1
$ws.Range("H130").Value2 = 50429
$ws.Range("J130").Value2 = 50429
$ws.Range("L130").Value2 = 50429
$ws.Range("N130").Value2 = -60469
$ws.Range("H132").Value2 = 7221.72
$ws.Range("I132").Value2 = 7012.5
$ws.Range("K132").Value2 = 21037.5
$ws.Range("M132").Value2 = -18507.5
$ws.Range("H133").Value2 = 57616.4
$ws.Range("J133").Value2 = 57616.4
$ws.Range("L133").Value2 = 57616.4
$ws.Range("N133").Value2 = -62676.4
$ws.Range("H136").Value2 = 4265.5454
$ws.Range("I136").Value2 = 1634.8182
$ws.Range("K136").Value2 = 4904.4546
$ws.Range("M136").Value2 = -2354.4546
$ws.Range("H140").Value2 = 77999.5
$ws.Range("J140").Value2 = 77999.5
$ws.Range("L140").Value2 = 77999.5
$ws.Range("N140").Value2 = -88359.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value2 = 21831.285
$ws.Range("I28").Value2 = 11950.5
$ws.Range("J28").Value2 = 25783.6
$ws.Range("K28").Value2 = 11950.5
$ws.Range("L28").Value2 = 25783.6
$ws.Range("M28").Value2 = -11602.5
$ws.Range("N28").Value2 = -26479.6
$ws.Range("H30").Value2 = 23850
$ws.Range("J30").Value2 = 30300
$ws.Range("L30").Value2 = 30300
$ws.Range("N30").Value2 = -30514
$ws.Range("H46").Value2 = 99996.5
$ws.Range("J46").Value2 = 99996.5
$ws.Range("L46").Value2 = 99996.5
$ws.Range("N46").Value2 = -100458.5
$ws.Range("H62").Value2 = 25044.666
$ws.Range("I62").Value2 = 45998.332
$ws.Range("J62").Value2 = 4091
$ws.Range("K62").Value2 = 45998.332
$ws.Range("L62").Value2 = 4091
$ws.Range("M62").Value2 = -45374.332
$ws.Range("N62").Value2 = -5339
$ws.Range("H63").Value2 = 31274.8
$ws.Range("J63").Value2 = 31274.8
$ws.Range("L63").Value2 = 31274.8
$ws.Range("N63").Value2 = -32522.8
$ws.Range("H65").Value2 = 25044.666
$ws.Range("I65").Value2 = 45998.332
$ws.Range("J65").Value2 = 4091
$ws.Range("K65").Value2 = 229991.66
$ws.Range("L65").Value2 = 20455
$ws.Range("M65").Value2 = -226871.66
$ws.Range("N65").Value2 = -26695
$ws.Range("H66").Value2 = 31274.8
$ws.Range("J66").Value2 = 31274.8
$ws.Range("L66").Value2 = 93824.4
$ws.Range("N66").Value2 = -100064.4
$ws.Range("H81").Value2 = 5784
$ws.Range("J81").Value2 = 9770
$ws.Range("L81").Value2 = 19540
$ws.Range("N81").Value2 = -21662
$ws.Range("H84").Value2 = 5784
$ws.Range("J84").Value2 = 9770
$ws.Range("L84").Value2 = 97700
$ws.Range("N84").Value2 = -108308
$ws.Range("H88").Value2 = 8750
$ws.Range("J88").Value2 = 8750
$ws.Range("L88").Value2 = 8750
$ws.Range("N88").Value2 = -9562
$ws.Range("H91").Value2 = 8750
$ws.Range("J91").Value2 = 8750
$ws.Range("L91").Value2 = 8750
$ws.Range("N91").Value2 = -11558
$ws.Range("H126").Value2 = 2630.7273
$ws.Range("I126").Value2 = 2630.7273
$ws.Range("J126").Value2 = 0
$ws.Range("K126").Value2 = 7892.1819
$ws.Range("L126").Value2 = 0
$ws.Range("M126").Value2 = -5422.1819
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value2 = 7150.4136
$ws.Range("I132").Value2 = 5141.8125
$ws.Range("K132").Value2 = 15425.4375
$ws.Range("M132").Value2 = -12895.4375
$ws.Range("H133").Value2 = 84499.75
$ws.Range("J133").Value2 = 84499.75
$ws.Range("L133").Value2 = 84499.75
$ws.Range("N133").Value2 = -94619.75
$ws.Range("H134").Value2 = 99996.5
$ws.Range("J134").Value2 = 99996.5
$ws.Range("L134").Value2 = 299989.5
$ws.Range("N134").Value2 = -305059.5
